$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Settings")

# Clear all existing hyperlinks first - they will be re-added at their final
# locations below (the engine only reliably supports clearing the whole
# collection, not individual items).
$ws.Hyperlinks.Delete()

# Row 7: strCashAccountsTemplateFilePath (was row 8) - keeps Hyperlink style
$ws.Range("A7").Value = "strCashAccountsTemplateFilePath"
$ws.Range("B7").Value = "\\LRRBTUIPFSP100\Profiles\Uipath_26\Desktop\EssilorLuxottica Projects\LUX-7.0_Running Bank Files in SAP\LUX - 7.0 - Retrieve SAP Bank Files\Data\"
$ws.Range("B7").Style = "Hyperlink"

# Row 8: strCashAccountsTemplateFileName (was row 9) - style becomes plain
$ws.Range("A8").Value = "strCashAccountsTemplateFileName"
$ws.Range("B8").Value = "CashAccountsTemplate.xlsx"
$ws.Range("B8").Style = "Normal"

# Row 9: strCashAccountsTemplateSheetName (was row 10) - stays plain
$ws.Range("A9").Value = "strCashAccountsTemplateSheetName"
$ws.Range("B9").Value = "Sheet1"
$ws.Range("B9").Style = "Normal"

# Row 10: boolAddHeaders (was row 11) - gets left/top alignment style
$ws.Range("A10").Value = "boolAddHeaders"
$ws.Range("B10").Value = $true
$ws.Range("B10").HorizontalAlignment = -4131
$ws.Range("B10").VerticalAlignment = -4160

# Row 11: strRange (was row 12) - style becomes plain
$ws.Range("A11").Value = "strRange"
$ws.Range("B11").Value = "A1"
$ws.Range("B11").Style = "Normal"

# Row 12: strSAPBankFilesPathPROD (NEW row) - Hyperlink style
$ws.Range("A12").Value = "strSAPBankFilesPathPROD"
$ws.Range("B12").Value = "\\myemdrive_server.luxgroup.net\MyEMDrive\EYEMED\ProductionData\SAP Auto Cash Application\Daily SAP Bank Files - ODS\"
$ws.Range("B12").Style = "Hyperlink"

# Row 13: strSAPBankFilesPathDEV (NEW row) - Hyperlink style
$ws.Range("A13").Value = "strSAPBankFilesPathDEV"
$ws.Range("B13").Value = "\\LRRBTUIPFSP100\Profiles\Uipath_26\Desktop\EssilorLuxottica Projects\LUX-7.0_Running Bank Files in SAP\LUX - 7.0 - Retrieve SAP Bank Files\Data\Input\"
$ws.Range("B13").Style = "Hyperlink"

# Re-create the hyperlinks at their final cell locations
$ws.Hyperlinks.Add($ws.Range("B6"), "../../../Components/")
$ws.Range("B6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B7"), ".")
$ws.Range("B7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B12"), "Input/")
$ws.Range("B12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B13"), "Input/")
$ws.Range("B13").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Assets sheet - new row referencing the SAP bank files path orchestrator asset
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Assets")
$ws3.Range("A3").Value = "strSAPBankFilesPath"
$ws3.Range("B3").Value = "Ess.LUX-7.0_SAPBankFilesPath"

# ---------------------------------------------------------------------------
# Window / selection state: Assets becomes the active tab with B5 selected,
# Settings keeps B17 selected (but is no longer the active tab).
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()
$ws3.Activate()
$ws3.Range("B5").Select()
